# agregue pame.html confirmacion de pedidos
#
# Updates the "Pedidos" sheet: row 2 and row 3 get new values (orders moved
# along in the workflow) and rows 4-8 are appended as new pedidos. Also adds
# a new user "Pame" to the "Usuarios" sheet (row 5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as genuine text, the same way Excel
# stores it when a user types a leading apostrophe in front of a number
# (keeps ID / phone-number-like strings from turning into numeric values).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = "'" + $Text
}

# ---------------------------------------------------------------------------
# Sheet "Pedidos"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pedidos")

# Row 2 - existing pedido updated further along the workflow
Set-TextValue $ws.Cells.Item(2, 1) "1739053700635"
$ws.Cells.Item(2, 2).Value = "Caipirinha"
$ws.Cells.Item(2, 3).Value = "1 litro(s)"
$ws.Cells.Item(2, 4).Value = "Transferencia"
$ws.Cells.Item(2, 5).Value = "manu"
Set-TextValue $ws.Cells.Item(2, 6) "595971224560"
$ws.Cells.Item(2, 7).Value = "8/2/2025, 19:28:24"
$ws.Cells.Item(2, 8).Value = "Terminado"
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 10).Value = "8/2/2025, 19:35:06"
$ws.Cells.Item(2, 11).Value = "Marcos"

# Row 3 - existing pedido updated further along the workflow
Set-TextValue $ws.Cells.Item(3, 1) "1739053945837"
$ws.Cells.Item(3, 2).Value = "Caipirinha"
$ws.Cells.Item(3, 3).Value = "1 litro(s)"
$ws.Cells.Item(3, 4).Value = "Transferencia"
$ws.Cells.Item(3, 5).Value = "Manu"
Set-TextValue $ws.Cells.Item(3, 6) "595971224560"
$ws.Cells.Item(3, 7).Value = "8/2/2025, 19:32:29"
$ws.Cells.Item(3, 8).Value = "Terminado"
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 10).Value = "8/2/2025, 19:35:07"
$ws.Cells.Item(3, 11).Value = "Marcos"

# Row 4 - new pedido (Pendiente)
Set-TextValue $ws.Cells.Item(4, 1) "1739054074695"
$ws.Cells.Item(4, 2).Value = "Caipirinha"
$ws.Cells.Item(4, 3).Value = "1 litro(s)"
$ws.Cells.Item(4, 4).Value = "Transferencia"
$ws.Cells.Item(4, 5).Value = "Chiqui Falcon"
Set-TextValue $ws.Cells.Item(4, 6) "595971224560"
$ws.Cells.Item(4, 7).Value = "8/2/2025, 19:34:43"
$ws.Cells.Item(4, 8).Value = "Pendiente"
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = "Pame"

# Row 5 - new pedido (Pendiente)
Set-TextValue $ws.Cells.Item(5, 1) "1739054236935"
$ws.Cells.Item(5, 2).Value = "Caipiruva"
$ws.Cells.Item(5, 3).Value = "2 litro(s)"
$ws.Cells.Item(5, 4).Value = "Efectivo"
$ws.Cells.Item(5, 5).Value = "gei"
Set-TextValue $ws.Cells.Item(5, 6) "595971224560"
$ws.Cells.Item(5, 7).Value = "8/2/2025, 19:37:18"
$ws.Cells.Item(5, 8).Value = "Pendiente"
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = "Pame"

# Row 6 - new pedido (Pendiente)
Set-TextValue $ws.Cells.Item(6, 1) "1739054255630"
$ws.Cells.Item(6, 2).Value = "Caipiruva"
$ws.Cells.Item(6, 3).Value = "1 litro(s)"
$ws.Cells.Item(6, 4).Value = "Transferencia"
$ws.Cells.Item(6, 5).Value = "manu"
Set-TextValue $ws.Cells.Item(6, 6) "595971224560"
$ws.Cells.Item(6, 7).Value = "8/2/2025, 19:37:41"
$ws.Cells.Item(6, 8).Value = "Pendiente"
$ws.Cells.Item(6, 9).Value = ""
$ws.Cells.Item(6, 10).Value = ""
$ws.Cells.Item(6, 11).Value = "Pame"

# Row 7 - new pedido (A Confirmar)
Set-TextValue $ws.Cells.Item(7, 1) "1739054326879"
$ws.Cells.Item(7, 2).Value = "Caipirinha"
$ws.Cells.Item(7, 3).Value = "1 litro(s)"
$ws.Cells.Item(7, 4).Value = "Transferencia"
$ws.Cells.Item(7, 5).Value = "Manuel Ayala"
Set-TextValue $ws.Cells.Item(7, 6) "595971224560"
$ws.Cells.Item(7, 7).Value = "8/2/2025, 19:38:55"
$ws.Cells.Item(7, 8).Value = "A Confirmar"
$ws.Cells.Item(7, 9).Value = ""
$ws.Cells.Item(7, 10).Value = ""
$ws.Cells.Item(7, 11).Value = ""

# Row 8 - new pedido (A Confirmar)
Set-TextValue $ws.Cells.Item(8, 1) "1739054446904"
$ws.Cells.Item(8, 2).Value = "Caipirinha"
$ws.Cells.Item(8, 3).Value = "1 litro(s)"
$ws.Cells.Item(8, 4).Value = "Efectivo"
$ws.Cells.Item(8, 5).Value = "Manuel Ayala"
Set-TextValue $ws.Cells.Item(8, 6) "595971224560"
$ws.Cells.Item(8, 7).Value = "8/2/2025, 19:40:57"
$ws.Cells.Item(8, 8).Value = "A Confirmar"
$ws.Cells.Item(8, 9).Value = ""
$ws.Cells.Item(8, 10).Value = ""
$ws.Cells.Item(8, 11).Value = ""

# ---------------------------------------------------------------------------
# Sheet "Usuarios" - add Pame as confirmation user
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Usuarios")
$ws2.Cells.Item(5, 1).Value = 4
$ws2.Cells.Item(5, 2).Value = "Pame"
$ws2.Cells.Item(5, 3).Value = 123456
